# The Dragon Boat Festival
#
# Paragraph 4 ("大晴") currently ends the document and carries the
# "_GoBack" bookmark plus a paragraph-mark font hint of "default".
# This edit:
#   1. changes paragraph 4's paragraph-mark font hint to "eastAsia"
#      and drops the "_GoBack" bookmark from it,
#   2. appends two new content paragraphs ("2022年6月22日星期三" and
#      "小雨" -- the latter now carrying the "_GoBack" bookmark), and
#   3. appends a final empty paragraph (paragraph-mark hint "default"),
# matching the tail of the document being rebuilt for a new diary entry.

$d = $word.ActiveDocument

$p4 = $d.Paragraphs.Item(4)
$p4Range = $p4.Range

$fragment = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>大晴</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2022年6月22日星期三</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>小雨</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$flatOpc = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML replaces the addressed range's contents, so targeting the whole
# of paragraph 4 (its run text AND its paragraph mark) with three <w:p>
# elements both rewrites paragraph 4 in place and appends the two new
# paragraphs after it; the paragraph mark that used to end paragraph 4
# survives as a new trailing empty paragraph (keeping its original
# rFonts hint="default"), which is exactly the blank paragraph the diff adds
# at the end of the document.
$p4Range.InsertXML($flatOpc)
